$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 353, shifting existing rows 353:426 down to 354:427
$ws.Rows.Item(353).Insert()

# Populate the new row 353 with the inserted record's data
$ws.Range("A353").Value = 9
$ws.Range("B353").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C353").Value = "Metropolitana"
$ws.Range("D353").Value = 45258
$ws.Range("E353").Value = 13
$ws.Range("F353").Value = 100112026
$ws.Range("G353").Value = "Haba"
$ws.Range("H353").Value = "Sin especificar"
$ws.Range("I353").Value = "Primera"
$ws.Range("J353").Value = 45
$ws.Range("K353").Value = 13000
$ws.Range("L353").Value = 15000
$ws.Range("M353").Value = 13889
$ws.Range("N353").Value = "$/saco 25 kilos"
$ws.Range("O353").Value = "Región del Maule"
$ws.Range("P353").Value = 556
$ws.Range("Q353").Value = 25
$ws.Range("R353").Value = "Hortaliza"

# Match the date-number formatting style used by the rest of column D
$ws.Range("D353").NumberFormat = $ws.Range("D354").NumberFormat
